$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for account 004472386 (GABRIEL) and 004547722 (MARCIA).
# These are currently Excel rows 4 and 5 (row 1 is the header).
$ws.Range("A4:A5").EntireRow.Delete()
